# This edit inserts one new data row (a new weekly observation) right after
# the existing row 152 (i.e. at row 153), shifting the previous rows 153-197
# down to 154-198, and populates the newly inserted row with its own data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 153, shifting existing rows 153:197 down to 154:198
$ws.Rows.Item(153).Insert()

# Populate the newly inserted row 153 with the new observation's data
$ws.Cells.Item(153, 1).Value  = 3
$ws.Cells.Item(153, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(153, 3).Value  = "Coquimbo"
$ws.Cells.Item(153, 4).Value  = 44463
$ws.Cells.Item(153, 5).Value  = 5
$ws.Cells.Item(153, 6).Value  = 100112009
$ws.Cells.Item(153, 7).Value  = "Acelga"
$ws.Cells.Item(153, 8).Value  = "Sin especificar"
$ws.Cells.Item(153, 9).Value  = "Primera"
$ws.Cells.Item(153, 10).Value = 240
$ws.Cells.Item(153, 11).Value = 2000
$ws.Cells.Item(153, 12).Value = 2200
$ws.Cells.Item(153, 13).Value = 2108
$ws.Cells.Item(153, 14).Value = "`$/docena de atados (6 kilos)"
$ws.Cells.Item(153, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(153, 16).Value = 351
$ws.Cells.Item(153, 17).Value = 6
$ws.Cells.Item(153, 18).Value = "Hortaliza"
